# manual-tests/shopping-test-suite.xlsx
# "added new sceanrio and tests"
#
# CHECK-002 (row 26) is renumbered to CHECK-003, and a brand-new CHECK-004
# scenario ("Check for required fields") is appended as row 27. The sheet's
# scroll position / selection is also updated to reflect where the author was
# working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite")
$ws.Activate()

# --- New row 27: CHECK-004 "Check for required fields" -----------------
# Values are written in the same first-use order the author's Excel session
# produced so newly interned shared strings line up the same way:
#   A27 (CHECK-004) -> A26 (CHECK-003) -> D27 -> E27 -> B27 -> C27 -> F27
$ws.Range("A27").Value = "CHECK-004"

# Existing CHECK-002 scenario becomes CHECK-003 (new test slotted in before it)
$ws.Range("A26").Value = "CHECK-003"

$ws.Range("D27").Value = "1) Users fills no details 2)  Click Place Order "
$ws.Range("E27").Value = "A message will say `"Full name is required. A valid email is required. Address is required. City is required. Postcode is required.`" "
$ws.Range("B27").Value = "Check for required fields  "
$ws.Range("C27").Value = "1) CART-001 completed 2) User clicks checkout "
$ws.Range("F27").Value = "P0"

# --- View state: scroll + selection -------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 7
$win.ScrollColumn = 5
$win.Top = 760
$win.Height = 17760

$null = $ws.Range("G23").Select()
